$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.975.46'
$ws.Range("E2").Value = '  +8.38%  '

$ws.Range("D3").Value = '3.437.73'
$ws.Range("E3").Value = '  +5.51%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '414.30'
$ws.Range("E5").Value = '  +4.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '123.82'
$ws.Range("E6").Value = '  +13.63%  '

$ws.Range("D7").Value = '3.432.96'
$ws.Range("E7").Value = '  +5.52%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.593'
$ws.Range("E8").Value = '  +1.87%  '

$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.657'
$ws.Range("E10").Value = '  +5.79%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.127'
$ws.Range("E11").Value = '  +32.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '41.34'
$ws.Range("E12").Value = '  +4.94%  '

$ws.Range("E13").Value = '  -0.30%  '

$ws.Range("D14").Value = '3.982.10'
$ws.Range("E14").Value = '  +5.52%  '

$ws.Range("E15").Value = '  +2.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.84'
$ws.Range("E16").Value = '  +4.38%  '

$ws.Range("D17").Value = '3.437.44'
$ws.Range("E17").Value = '  +5.42%  '

$ws.Range("D18").Value = '61.963.73'
$ws.Range("E18").Value = '  +8.67%  '

$ws.Range("E19").Value = '  -0.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.96'
$ws.Range("E20").Value = '  -1.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000131'
$ws.Range("E21").Value = '  +20.52%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.35'
$ws.Range("E22").Value = '  +0.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '82.35'
$ws.Range("E23").Value = '  +10.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '314.29'
$ws.Range("E24").Value = '  +6.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.98'
$ws.Range("E25").Value = '  +0.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.17'
$ws.Range("E26").Value = '  -0.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '31.00'
$ws.Range("E27").Value = '  +10.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.84'
$ws.Range("E28").Value = '  +5.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.87'
$ws.Range("E29").Value = '  -0.62%  '

$ws.Range("E30").Value = '  -2.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.173'
$ws.Range("E31").Value = '  +2.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.117'
$ws.Range("E32").Value = '  +4.80%  '

$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.57'
$ws.Range("E33").Value = '  +3.25%  '

$ws.Range("B34").Value = 'Toncoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.58'
$ws.Range("E34").Value = '  +20.60%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '42.28'
$ws.Range("E35").Value = '  +5.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0485'
$ws.Range("E37").Value = '  -0.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.27'
$ws.Range("E38").Value = '  +1.87%  '

$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.51'
$ws.Range("E39").Value = '  +1.41%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.00'
$ws.Range("E41").Value = '  -0.25%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.99'
$ws.Range("E42").Value = '  +6.47%  '

$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.125'
$ws.Range("E43").Value = '  +3.46%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '133.94'
$ws.Range("E44").Value = '  -1.87%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.23'
$ws.Range("E45").Value = '  +2.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.284'
$ws.Range("E46").Value = '  +0.29%  '

$ws.Range("E47").Value = '  -1.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.07'
$ws.Range("E48").Value = '  -1.94%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.20'
$ws.Range("E49").Value = '  -0.30%  '

$ws.Range("D50").Value = '2.205.49'
$ws.Range("E50").Value = '  +2.67%  '

$ws.Range("D51").Value = '3.779.32'
$ws.Range("E51").Value = '  +5.42%  '
